# Auto-generated edit script: applies per-cell numeric updates
# matching the target diff for Sheets/Tonberry_Profits.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 408.33334
$ws.Range("I2").Value = 233.33333
$ws.Range("K2").Value = 233.33333
$ws.Range("M2").Value = -120.33333

$ws.Range("H62").Value = 2829.2856
$ws.Range("I62").Value = 2701.25
$ws.Range("K62").Value = 2701.25
$ws.Range("M62").Value = -2077.25

$ws.Range("H65").Value = 2829.2856
$ws.Range("I65").Value = 2701.25
$ws.Range("K65").Value = 13506.25
$ws.Range("M65").Value = -10386.25

$ws.Range("H132").Value = 982.3396
$ws.Range("I132").Value = 829.7755
$ws.Range("K132").Value = 2489.3265
$ws.Range("M132").Value = 40.67349999999988

$ws.Range("H135").Value = 603.5
$ws.Range("I135").Value = 515.44446
$ws.Range("K135").Value = 4639.00014
$ws.Range("M135").Value = -2104.00014

$ws.Range("H137").Value = 1573.3158
$ws.Range("I137").Value = 1361.8462
$ws.Range("J137").Value = 2031.5
$ws.Range("K137").Value = 4085.5386
$ws.Range("L137").Value = 6094.5
$ws.Range("M137").Value = -1535.5386
$ws.Range("N137").Value = -11194.5

$ws.Range("H138").Value = 3482.8
$ws.Range("J138").Value = 3850
$ws.Range("L138").Value = 11550
$ws.Range("N138").Value = -21830

$ws.Range("H141").Value = 802365
$ws.Range("I141").Value = 967176.4399999999
$ws.Range("K141").Value = 2901529.32
$ws.Range("M141").Value = -2896349.32

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2567.5325
$ws.Range("I32").Value = 1936.8955
$ws.Range("K32").Value = 1936.8955
$ws.Range("M32").Value = -1649.8955

$ws.Range("H109").Value = 47184
$ws.Range("J109").Value = 47184
$ws.Range("L109").Value = 47184
$ws.Range("N109").Value = -49958

$ws.Range("H132").Value = 1827.1111
$ws.Range("I132").Value = 1214.3529
$ws.Range("K132").Value = 3643.0587
$ws.Range("M132").Value = -1113.0587

$ws.Range("H139").Value = 40150.75
$ws.Range("J139").Value = 40150.75
$ws.Range("L139").Value = 40150.75
$ws.Range("N139").Value = -50430.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2063.0667
$ws.Range("J20").Value = 2416.25
$ws.Range("L20").Value = 2416.25
$ws.Range("N20").Value = -2910.25

$ws.Range("H75").Value = 9252.333000000001
$ws.Range("I75").Value = 7945
$ws.Range("K75").Value = 7945
$ws.Range("M75").Value = -7009

$ws.Range("H78").Value = 9252.333000000001
$ws.Range("I78").Value = 7945
$ws.Range("K78").Value = 23835
$ws.Range("M78").Value = -19155

$ws.Range("H99").Value = 1592.5
$ws.Range("I99").Value = 1345.5555
$ws.Range("K99").Value = 1345.5555
$ws.Range("M99").Value = 152.4445000000001

$ws.Range("H134").Value = 3682.4
$ws.Range("I134").Value = 3985.3408
$ws.Range("J134").Value = 1460.8334
$ws.Range("K134").Value = 11956.0224
$ws.Range("L134").Value = 4382.5002
$ws.Range("M134").Value = -9421.0224
$ws.Range("N134").Value = -9452.5002

$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 823
$ws.Range("I16").Value = 775.5
$ws.Range("K16").Value = 775.5
$ws.Range("M16").Value = -488.5

$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").Value = $null

$ws.Range("H86").Value = 2099.3
$ws.Range("I86").Value = 1699.2858
$ws.Range("J86").Value = 3032.6667
$ws.Range("K86").Value = 1699.2858
$ws.Range("L86").Value = 3032.6667
$ws.Range("M86").Value = -576.2858000000001
$ws.Range("N86").Value = -5278.6667

$ws.Range("H89").Value = 2099.3
$ws.Range("I89").Value = 1699.2858
$ws.Range("J89").Value = 3032.6667
$ws.Range("K89").Value = 8496.429
$ws.Range("L89").Value = 15163.3335
$ws.Range("M89").Value = -2880.429
$ws.Range("N89").Value = -26395.3335

$ws.Range("H105").Value = 1935
$ws.Range("I105").Value = 1822
$ws.Range("K105").Value = 1822
$ws.Range("M105").Value = -75

$ws.Range("H107").Value = 397.64706
$ws.Range("I107").Value = 321.42856
$ws.Range("J107").Value = 753.3333
$ws.Range("K107").Value = 321.42856
$ws.Range("L107").Value = 753.3333
$ws.Range("M107").Value = 1598.57144
$ws.Range("N107").Value = -4593.3333

$ws.Range("H113").Value = 823
$ws.Range("I113").Value = 775.5
$ws.Range("K113").Value = 775.5
$ws.Range("M113").Value = 1394.5

$ws.Range("H122").Value = 6539.6
$ws.Range("I122").Value = 9000
$ws.Range("K122").Value = 27000
$ws.Range("M122").Value = -24550

$ws.Range("H131").Value = 38883.418
$ws.Range("J131").Value = 38883.418
$ws.Range("L131").Value = 38883.418
$ws.Range("N131").Value = -48963.418

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 773.4286
$ws.Range("I26").Value = 927.25
$ws.Range("K26").Value = 2781.75
$ws.Range("M26").Value = -2493.75

$ws.Range("H107").Value = 758.1905
$ws.Range("J107").Value = 935.73334
$ws.Range("L107").Value = 2807.20002
$ws.Range("N107").Value = -6647.20002

$ws.Range("H122").Value = 1012.1
$ws.Range("I122").Value = 585
$ws.Range("J122").Value = 1296.8334
$ws.Range("K122").Value = 5265
$ws.Range("L122").Value = 11671.5006
$ws.Range("M122").Value = -2815
$ws.Range("N122").Value = -16571.5006

$ws.Range("H131").Value = 8718.204
$ws.Range("J131").Value = 9655.379999999999
$ws.Range("L131").Value = 28966.14
$ws.Range("N131").Value = -39046.14

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 17987.25
$ws.Range("J92").Value = 17987.25
$ws.Range("L92").Value = 17987.25
$ws.Range("N92").Value = -21731.25

$ws.Range("H122").Value = 1571.129
$ws.Range("J122").Value = 2005.3334
$ws.Range("L122").Value = 6016.0002
$ws.Range("N122").Value = -10916.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1926.3636
$ws.Range("I22").Value = 1961.25
$ws.Range("J22").Value = 1833.3334
$ws.Range("K22").Value = 1961.25
$ws.Range("L22").Value = 1833.3334
$ws.Range("M22").Value = -1666.25
$ws.Range("N22").Value = -2423.3334

$ws.Range("H27").Value = 1926.3636
$ws.Range("I27").Value = 1961.25
$ws.Range("J27").Value = 1833.3334
$ws.Range("K27").Value = 1961.25
$ws.Range("L27").Value = 1833.3334
$ws.Range("M27").Value = -1854.25
$ws.Range("N27").Value = -2047.3334

$ws.Range("H61").Value = 2977.4443
$ws.Range("I61").Value = 2899.5715
$ws.Range("K61").Value = 2899.5715
$ws.Range("M61").Value = -2697.5715

$ws.Range("H100").Value = 1493.25
$ws.Range("I100").Value = 1493.25
$ws.Range("K100").Value = 1493.25
$ws.Range("M100").Value = -952.25

$ws.Range("H113").Value = 2977.4443
$ws.Range("I113").Value = 2899.5715
$ws.Range("K113").Value = 2899.5715
$ws.Range("M113").Value = -729.5715

$ws.Range("H122").Value = 6685.3335
$ws.Range("J122").Value = 11666.667
$ws.Range("L122").Value = 35000.001
$ws.Range("N122").Value = -39900.001

$ws.Range("H132").Value = 1835.919
$ws.Range("I132").Value = 1607.5555
$ws.Range("J132").Value = 2052.2632
$ws.Range("K132").Value = 4822.666499999999
$ws.Range("L132").Value = 6156.7896
$ws.Range("M132").Value = -2292.666499999999
$ws.Range("N132").Value = -11216.7896

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 47343.5
$ws.Range("J95").Value = 47343.5
$ws.Range("L95").Value = 47343.5
$ws.Range("N95").Value = -52835.5

$ws.Range("H107").Value = 831.36365
$ws.Range("I107").Value = 642.75
$ws.Range("K107").Value = 1928.25
$ws.Range("M107").Value = -8.25

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = $null

$ws.Range("H119").Value = 28694
$ws.Range("J119").Value = 28694
$ws.Range("L119").Value = 28694
$ws.Range("N119").Value = -38370

$ws.Range("H122").Value = 71943.45
$ws.Range("I122").Value = 71943.45
$ws.Range("K122").Value = 215830.35
$ws.Range("M122").Value = -213380.35

$ws.Range("H132").Value = 1168.9744
$ws.Range("I132").Value = 869.4
$ws.Range("K132").Value = 2608.2
$ws.Range("M132").Value = -78.19999999999982
